$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")
$shp = $ws.Shapes.Item(1)
$tf2 = $shp.TextFrame2
$tr = $tf2.TextRange
$paras = $tr.Paragraphs
Write-Host "count before" $paras.Count
$paras.Item(1).Text = "FIRST PARA ONLY"
